$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the status value up from C4 to C2, with updated text ("DONE" -> "Completed"),
# and remove the now-empty row 4.
$ws.Range("C2").Value = "Completed"
$ws.Range("C4").Value = $null
$ws.Rows.Item(4).Delete()
